# Fix bug with column colors (Z, AA, AB) on rows 2-8: the zebra-stripe
# fill pattern (alternating white / grey) was broken for these three
# columns; restore the correct alternating colors.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$white = 16777215   # RGB(255,255,255)
$grey  = 14540253   # RGB(221,221,221) / 0xDDDDDD

foreach ($row in 2..8) {
    $ws.Range("Z$row").Interior.Color  = $grey
    $ws.Range("AA$row").Interior.Color = $white
    $ws.Range("AB$row").Interior.Color = $grey
}

# Adjust page margins back to Excel's defaults (inches): 0.75/0.75/1/1/0.5/0.5
# PageSetup margins are expressed in points (1 inch = 72 points).
$ws.PageSetup.LeftMargin   = 54   # 0.75 in
$ws.PageSetup.RightMargin  = 54   # 0.75 in
$ws.PageSetup.TopMargin    = 72   # 1 in
$ws.PageSetup.BottomMargin = 72   # 1 in
$ws.PageSetup.HeaderMargin = 36   # 0.5 in
$ws.PageSetup.FooterMargin = 36   # 0.5 in

# Add workbook protection element (workbookProtection) as in the target file.
$wb.Protect()
